$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 4333.3335

$ws.Range("H103").Value = 1215.7693
$ws.Range("I103").Value = 1171.4286
$ws.Range("J103").Value = 1267.5
$ws.Range("K103").Value = 3514.2858
$ws.Range("L103").Value = 3802.5
$ws.Range("M103").Value = -2928.2858
$ws.Range("N103").Value = -4974.5

$ws.Range("H113").Value = 1961.9166
$ws.Range("I113").Value = 2448.3333
$ws.Range("J113").Value = 1799.7778
$ws.Range("K113").Value = 2448.3333
$ws.Range("L113").Value = 1799.7778
$ws.Range("M113").Value = 805.6667000000002
$ws.Range("N113").Value = -8307.7778

$ws.Range("H134").Value = 191233.33
$ws.Range("J134").Value = 191233.33
$ws.Range("L134").Value = 191233.33
$ws.Range("N134").Value = -201373.33

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5507.0713
$ws.Range("I74").Value = 1260.5714
$ws.Range("K74").Value = 1260.5714
$ws.Range("M74").Value = -386.5714

$ws.Range("H77").Value = 5507.0713
$ws.Range("I77").Value = 1260.5714
$ws.Range("K77").Value = 6302.857
$ws.Range("M77").Value = -1934.857

$ws.Range("H132").Value = 1453.4595
$ws.Range("I132").Value = 1236.875
$ws.Range("K132").Value = 3710.625
$ws.Range("M132").Value = -1180.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H64").Value = 1641.96
$ws.Range("I64").Value = 3604
$ws.Range("J64").Value = 878.94446
$ws.Range("K64").Value = 3604
$ws.Range("L64").Value = 878.94446
$ws.Range("M64").Value = -3379
$ws.Range("N64").Value = -1328.94446

$ws.Range("H67").Value = 1641.96
$ws.Range("I67").Value = 3604
$ws.Range("J67").Value = 878.94446
$ws.Range("K67").Value = 3604
$ws.Range("L67").Value = 878.94446
$ws.Range("M67").Value = -2824
$ws.Range("N67").Value = -2438.94446

$ws.Range("H105").Value = 4300
$ws.Range("I105").Value = 10000
$ws.Range("K105").Value = 10000
$ws.Range("M105").Value = -8253

$ws.Range("H134").Value = 39023.566
$ws.Range("I134").Value = 41274.75
$ws.Range("K134").Value = 123824.25
$ws.Range("M134").Value = -121289.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1418.5
$ws.Range("I58").Value = 1153.0938
$ws.Range("J58").Value = 2267.8
$ws.Range("K58").Value = 1153.0938
$ws.Range("L58").Value = 2267.8
$ws.Range("M58").Value = -950.0938000000001
$ws.Range("N58").Value = -2673.8

$ws.Range("H62").Value = 2933.3333
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 2933.3333
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240

$ws.Range("H75").Value = 20999
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 20999
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 20999
$ws.Range("N75").Value = -22995
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 20999
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 20999
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 62997
$ws.Range("N78").Value = -72981
$ws.Range("M78").ClearContents()

$ws.Range("H132").Value = 1826.8055
$ws.Range("I132").Value = 1135
$ws.Range("K132").Value = 3405
$ws.Range("M132").Value = -875

$ws.Range("H134").Value = 2008.8077
$ws.Range("I134").Value = 1874.9546
$ws.Range("K134").Value = 5624.8638
$ws.Range("M134").Value = -3089.8638

$ws.Range("H136").Value = 1418.5
$ws.Range("I136").Value = 1153.0938
$ws.Range("J136").Value = 2267.8
$ws.Range("K136").Value = 3459.2814
$ws.Range("L136").Value = 6803.400000000001
$ws.Range("M136").Value = -909.2814000000003
$ws.Range("N136").Value = -11903.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1099.875
$ws.Range("I68").Value = 1150
$ws.Range("J68").Value = 1083.1666
$ws.Range("K68").Value = 3450
$ws.Range("L68").Value = 3249.4998
$ws.Range("M68").Value = -2639
$ws.Range("N68").Value = -4871.4998

$ws.Range("H71").Value = 1099.875
$ws.Range("I71").Value = 1150
$ws.Range("J71").Value = 1083.1666
$ws.Range("K71").Value = 10350
$ws.Range("L71").Value = 9748.499400000001
$ws.Range("M71").Value = -6294
$ws.Range("N71").Value = -17860.4994

$ws.Range("H76").Value = 3895
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 3895
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 35000
$ws.Range("J108").Value = 35000
$ws.Range("L108").Value = 35000
$ws.Range("N108").Value = -42680

$ws.Range("H122").Value = 4286.25
$ws.Range("I122").Value = 4578.913
$ws.Range("K122").Value = 13736.739
$ws.Range("M122").Value = -11286.739

$ws.Range("H132").Value = 3430.7827
$ws.Range("I132").Value = 3253.75
$ws.Range("K132").Value = 9761.25
$ws.Range("M132").Value = -7231.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1891.909
$ws.Range("I61").Value = 1038.875
$ws.Range("K61").Value = 1038.875
$ws.Range("M61").Value = -836.875

$ws.Range("H93").Value = 2399.8333
$ws.Range("I93").Value = 2480
$ws.Range("K93").Value = 2480
$ws.Range("M93").Value = -1232

$ws.Range("H113").Value = 1891.909
$ws.Range("I113").Value = 1038.875
$ws.Range("K113").Value = 1038.875
$ws.Range("M113").Value = 1131.125

$ws.Range("H122").Value = 4791.8184
$ws.Range("I122").Value = 5227.4136
$ws.Range("J122").Value = 3949.6667
$ws.Range("K122").Value = 15682.2408
$ws.Range("L122").Value = 11849.0001
$ws.Range("M122").Value = -13232.2408
$ws.Range("N122").Value = -16749.0001

$ws.Range("H132").Value = 6741.75
$ws.Range("I132").Value = 8984.923000000001
$ws.Range("K132").Value = 26954.769
$ws.Range("M132").Value = -24424.769

$ws.Range("H136").Value = 1805.6666
$ws.Range("I136").Value = 766.8
$ws.Range("K136").Value = 2300.4
$ws.Range("M136").Value = 249.6000000000004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 22500
$ws.Range("J28").Value = 32000
$ws.Range("L28").Value = 32000
$ws.Range("N28").Value = -32696

$ws.Range("H107").Value = 547.8570999999999
$ws.Range("I107").Value = 352.07144
$ws.Range("J107").Value = 939.4286
$ws.Range("K107").Value = 1056.21432
$ws.Range("L107").Value = 2818.2858
$ws.Range("M107").Value = 863.78568
$ws.Range("N107").Value = -6658.2858

$ws.Range("H132").Value = 1403.6471
$ws.Range("I132").Value = 808.84
$ws.Range("J132").Value = 3055.889
$ws.Range("K132").Value = 2426.52
$ws.Range("L132").Value = 9167.667000000001
$ws.Range("M132").Value = 103.48
$ws.Range("N132").Value = -14227.667

$ws.Range("H136").Value = 8522.593000000001
$ws.Range("I136").Value = 8804.4
$ws.Range("K136").Value = 26413.2
$ws.Range("M136").Value = -23863.2
